$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2433.3333
$ws.Range("I4").Value = 2150
$ws.Range("K4").Value = 2150
$ws.Range("M4").Value = -2036
$ws.Range("H33").Value = 130.1875
$ws.Range("I33").Value = 160.16667
$ws.Range("J33").Value = 112.2
$ws.Range("K33").Value = 160.16667
$ws.Range("L33").Value = 112.2
$ws.Range("M33").Value = 68.83332999999999
$ws.Range("N33").Value = -570.2
$ws.Range("H80").Value = 2960.25
$ws.Range("I80").Value = 3437
$ws.Range("J80").Value = 2165.6667
$ws.Range("K80").Value = 10311
$ws.Range("L80").Value = 6497.000100000001
$ws.Range("M80").Value = -9313
$ws.Range("N80").Value = -8493.000100000001
$ws.Range("H83").Value = 2960.25
$ws.Range("I83").Value = 3437
$ws.Range("J83").Value = 2165.6667
$ws.Range("K83").Value = 30933
$ws.Range("L83").Value = 19491.0003
$ws.Range("M83").Value = -25941
$ws.Range("N83").Value = -29475.0003
$ws.Range("H97").Value = 1221.75
$ws.Range("J97").Value = 1333
$ws.Range("L97").Value = 3999
$ws.Range("N97").Value = -4991
$ws.Range("H98").Value = 821.32355
$ws.Range("I98").Value = 536.1923
$ws.Range("K98").Value = 536.1923
$ws.Range("M98").Value = 961.8077
$ws.Range("H112").Value = 1366.5
$ws.Range("J112").Value = 1366.5
$ws.Range("L112").Value = 4099.5
$ws.Range("N112").Value = -6315.5
$ws.Range("H122").Value = 821.32355
$ws.Range("I122").Value = 536.1923
$ws.Range("K122").Value = 1608.5769
$ws.Range("M122").Value = 841.4231
$ws.Range("H125").Value = 365.3846
$ws.Range("I125").Value = 370.83334
$ws.Range("K125").Value = 3337.50006
$ws.Range("M125").Value = -877.5000600000003
$ws.Range("H137").Value = 2319.9333
$ws.Range("I137").Value = 1782.9
$ws.Range("J137").Value = 3394
$ws.Range("K137").Value = 5348.700000000001
$ws.Range("L137").Value = 10182
$ws.Range("M137").Value = -2798.700000000001
$ws.Range("N137").Value = -15282
$ws.Range("H138").Value = 1873.8108
$ws.Range("I138").Value = 1681.2222
$ws.Range("J138").Value = 2056.2632
$ws.Range("K138").Value = 5043.6666
$ws.Range("L138").Value = 6168.7896
$ws.Range("M138").Value = 96.33340000000044
$ws.Range("N138").Value = -16448.7896

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1492.6666
$ws.Range("I122").Value = 989
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 2967
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -517
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 1441.7142
$ws.Range("I132").Value = 1346.4166
$ws.Range("J132").Value = 2013.5
$ws.Range("K132").Value = 4039.2498
$ws.Range("L132").Value = 6040.5
$ws.Range("M132").Value = -1509.2498
$ws.Range("N132").Value = -11100.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 633798.7
$ws.Range("I86").Value = 772088.9
$ws.Range("J86").Value = 334170
$ws.Range("K86").Value = 772088.9
$ws.Range("L86").Value = 334170
$ws.Range("M86").Value = -770965.9
$ws.Range("N86").Value = -336416
$ws.Range("H89").Value = 633798.7
$ws.Range("I89").Value = 772088.9
$ws.Range("J89").Value = 334170
$ws.Range("K89").Value = 3860444.5
$ws.Range("L89").Value = 1670850
$ws.Range("M89").Value = -3854828.5
$ws.Range("N89").Value = -1682082
$ws.Range("H99").Value = 1089.7778
$ws.Range("J99").Value = 1173.5714
$ws.Range("L99").Value = 1173.5714
$ws.Range("N99").Value = -4169.5714
$ws.Range("H134").Value = 7966.8125
$ws.Range("I134").Value = 9818.583000000001
$ws.Range("J134").Value = 2411.5
$ws.Range("K134").Value = 29455.749
$ws.Range("L134").Value = 7234.5
$ws.Range("M134").Value = -26920.749
$ws.Range("N134").Value = -12304.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 752.3333
$ws.Range("I19").Value = 752.3333
$ws.Range("K19").Value = 752.3333
$ws.Range("M19").Value = -582.3333
$ws.Range("H22").Value = 1339.2
$ws.Range("I22").Value = 600.3333
$ws.Range("J22").Value = 1523.9166
$ws.Range("K22").Value = 600.3333
$ws.Range("L22").Value = 1523.9166
$ws.Range("M22").Value = -250.3333
$ws.Range("N22").Value = -2223.9166
$ws.Range("H24").Value = 752.3333
$ws.Range("I24").Value = 752.3333
$ws.Range("K24").Value = 752.3333
$ws.Range("M24").Value = -582.3333
$ws.Range("H31").Value = 2959.75
$ws.Range("I31").Value = 2726
$ws.Range("J31").Value = 3427.25
$ws.Range("K31").Value = 2726
$ws.Range("L31").Value = 3427.25
$ws.Range("M31").Value = -2431
$ws.Range("N31").Value = -4017.25
$ws.Range("H34").Value = 2959.75
$ws.Range("I34").Value = 2726
$ws.Range("J34").Value = 3427.25
$ws.Range("K34").Value = 2726
$ws.Range("L34").Value = 3427.25
$ws.Range("M34").Value = -2524
$ws.Range("N34").Value = -3831.25
$ws.Range("H58").Value = 5437887.5
$ws.Range("I58").Value = 21739530
$ws.Range("K58").Value = 21739530
$ws.Range("M58").Value = -21739327
$ws.Range("H62").Value = 2488.5
$ws.Range("I62").Value = 2318.3333
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 2318.3333
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -1694.3333
$ws.Range("N62").Value = -4247
$ws.Range("H65").Value = 2488.5
$ws.Range("I65").Value = 2318.3333
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 11591.6665
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -8471.666499999999
$ws.Range("N65").Value = -21235
$ws.Range("H103").Value = 17120
$ws.Range("I103").Value = 16255.75
$ws.Range("K103").Value = 16255.75
$ws.Range("M103").Value = -15083.75
$ws.Range("H122").Value = 4493.923
$ws.Range("I122").Value = 3582.6365
$ws.Range("J122").Value = 9506
$ws.Range("K122").Value = 10747.9095
$ws.Range("L122").Value = 28518
$ws.Range("M122").Value = -8297.9095
$ws.Range("N122").Value = -33418
$ws.Range("H132").Value = 1867.6562
$ws.Range("I132").Value = 1393.5652
$ws.Range("J132").Value = 3079.2222
$ws.Range("K132").Value = 4180.6956
$ws.Range("L132").Value = 9237.6666
$ws.Range("M132").Value = -1650.6956
$ws.Range("N132").Value = -14297.6666
$ws.Range("H136").Value = 5437887.5
$ws.Range("I136").Value = 21739530
$ws.Range("K136").Value = 65218590
$ws.Range("M136").Value = -65216040

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1070
$ws.Range("J11").Value = 1484.5
$ws.Range("L11").Value = 4453.5
$ws.Range("N11").Value = -4733.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 900
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 900
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 98
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 900
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 492
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 2119.8
$ws.Range("I122").Value = 1899.5
$ws.Range("J122").Value = 2266.6667
$ws.Range("K122").Value = 5698.5
$ws.Range("L122").Value = 6800.000100000001
$ws.Range("M122").Value = -3248.5
$ws.Range("N122").Value = -11700.0001
$ws.Range("H126").Value = 1952791.1
$ws.Range("I126").Value = 7939495.5
$ws.Range("K126").Value = 23818486.5
$ws.Range("M126").Value = -23816016.5
$ws.Range("H132").Value = 2568334.8
$ws.Range("I132").Value = 3500156.8
$ws.Range("K132").Value = 10500470.4
$ws.Range("M132").Value = -10497940.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1912.2307
$ws.Range("I7").Value = 1662.421
$ws.Range("J7").Value = 2590.2856
$ws.Range("K7").Value = 1662.421
$ws.Range("L7").Value = 2590.2856
$ws.Range("M7").Value = -1550.421
$ws.Range("N7").Value = -2814.2856
$ws.Range("H40").Value = 8922.809999999999
$ws.Range("I40").Value = 10037.308
$ws.Range("J40").Value = 7111.75
$ws.Range("K40").Value = 10037.308
$ws.Range("L40").Value = 7111.75
$ws.Range("M40").Value = -9901.308000000001
$ws.Range("N40").Value = -7383.75
$ws.Range("H82").Value = 1936.25
$ws.Range("I82").Value = 1307.2222
$ws.Range("J82").Value = 3823.3333
$ws.Range("K82").Value = 1307.2222
$ws.Range("L82").Value = 3823.3333
$ws.Range("M82").Value = -946.2221999999999
$ws.Range("N82").Value = -4545.3333
$ws.Range("H85").Value = 1936.25
$ws.Range("I85").Value = 1307.2222
$ws.Range("J85").Value = 3823.3333
$ws.Range("K85").Value = 1307.2222
$ws.Range("L85").Value = 3823.3333
$ws.Range("M85").Value = -59.22219999999993
$ws.Range("N85").Value = -6319.3333
$ws.Range("H122").Value = 14949.75
$ws.Range("I122").Value = 12899.5
$ws.Range("J122").Value = 17000
$ws.Range("K122").Value = 38698.5
$ws.Range("L122").Value = 51000
$ws.Range("M122").Value = -36248.5
$ws.Range("N122").Value = -55900
$ws.Range("H126").Value = 1912.2307
$ws.Range("I126").Value = 1662.421
$ws.Range("J126").Value = 2590.2856
$ws.Range("K126").Value = 4987.263
$ws.Range("L126").Value = 7770.8568
$ws.Range("M126").Value = -2517.263
$ws.Range("N126").Value = -12710.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 4180
$ws.Range("J62").Value = 1400
$ws.Range("K62").Value = 4180
$ws.Range("L62").Value = 1400
$ws.Range("M62").Value = -3556
$ws.Range("N62").Value = -2648
$ws.Range("I65").Value = 4180
$ws.Range("J65").Value = 1400
$ws.Range("K65").Value = 20900
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = -17780
$ws.Range("N65").Value = -13240
$ws.Range("H100").Value = 1075.5
$ws.Range("I100").Value = 930.6
$ws.Range("K100").Value = 1861.2
$ws.Range("M100").Value = -1320.2
$ws.Range("H126").Value = 1893.0435
$ws.Range("I126").Value = 1599
$ws.Range("K126").Value = 4797
$ws.Range("M126").Value = -2327
$ws.Range("H136").Value = 15874954
$ws.Range("I136").Value = 26456966
$ws.Range("J136").Value = 1936.4286
$ws.Range("K136").Value = 79370898
$ws.Range("L136").Value = 5809.2858
$ws.Range("M136").Value = -79368348
$ws.Range("N136").Value = -10909.2858
